$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column H to fit the longer "CodeQL/Bandit" detector label
$ws.Columns.Item(8).ColumnWidth = 24.8

# Update detector labels for CWE-022 and CWE-078 rows
$ws.Range("H3").Value = "CodeQL/Bandit"
$ws.Range("H4").Value = "CodeQL/Bandit"

# CWE-079 row: one more analyzer hit
$ws.Range("B5").Value = 2

# CWE-080 row: swap the automated hit for a manually-found one
$ws.Range("B6").Value = 0
$ws.Range("E6").Value = 1
$ws.Range("H6").Value = "Manual"

# CWE-089 row: add an Authors hit
$ws.Range("E7").Value = 1

# Move the active selection to E7, matching the author's last edit position
$ws.Range("E7").Select()
